{"js": "// Replace the four text passages described by the diff. Each old string is\n// the full, exact text of a single run's <w:t>, so a straightforward\n// search + replace (matching the whole paragraph text, case-sensitive) is\n// safe and will not clobber sibling runs or formatting.\n//\n// NOTE: the title \"Creating Stunning Digital Graphics with Canva\" is itself\n// a substring of the longer intro-sentence run, so the longer passage MUST\n// be replaced first \u2014 otherwise replacing the short title first would also\n// rewrite the matching substring inside the sentence, and the sentence's\n// old-text search would then no longer find a match.\nconst replacements = [\n  {\n    oldText:\n      \"We are applying for WSQ funding support for this new course Creating Stunning Digital Graphics with Canva according to Visual Communication DSN-COM-4005-1.1 under Design Framework.\",\n    newText:\n      \"We are applying for WSQ funding support for this new course AZ-900 Microsoft Azure Fundamentals Certification according to Cloud Computing ICT-DIT-4020-1.1 under Infocomm Technology Framework.\"\n  },\n  {\n    oldText: \"Creating Stunning Digital Graphics with Canva\",\n    newText: \"AZ-900 Microsoft Azure Fundamentals Certification\"\n  },\n  {\n    oldText:\n      \"A significant challenge in the field is the inability to create coherent storyboards that clearly communicate visual intentions. Many designers lack the skills to evaluate and enhance their designs, leading to missed opportunities for improvement. This gap can hinder collaboration among team members and affect the overall quality of visual projects.\",\n    newText:\n      \"A prevalent issue in the industry is the inadequate integration of existing systems with cloud computing components. This often leads to compatibility problems and hinders the full utilization of cloud capabilities. Furthermore, organizations frequently lack the necessary skills to conduct impact analysis, which is crucial for understanding the implications of cloud adoption on their operations.\"\n  },\n  {\n    oldText:\n      \"By emphasizing the creation of storyboards and task flows, this course addresses the need for clear visual communication. Participants will learn to evaluate their designs critically and suggest enhancements, fostering a collaborative environment that improves project outcomes.\",\n    newText:\n      \"The course equips participants with the skills to draft specifications that address performance and security requirements, ensuring that cloud solutions are tailored to organizational needs. Additionally, it emphasizes the importance of aligning existing systems with cloud components, which can significantly improve integration and overall system performance.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the four text passages described by the diff.\n#\n# Each old string is the full, exact text content of its paragraph (the\n# paragraph's only text-bearing run; any sibling runs are empty formatting\n# placeholders). Matching on the whole paragraph text and assigning\n# `Paragraph.Range.Text` in place swaps just that run's <w:t> and leaves\n# every other run/paragraph (including the empty sibling runs) untouched.\n#\n# (Using Find/Replace with ReplaceAll across the whole document instead\n# would also locate the text, but when the match length changes it can\n# collapse/merge neighbouring empty runs in the same paragraph as a side\n# effect \u2014 so the more surgical per-paragraph text assignment below is used\n# instead.)\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{\n    Old = \"Creating Stunning Digital Graphics with Canva\"\n    New = \"AZ-900 Microsoft Azure Fundamentals Certification\"\n  },\n  @{\n    Old = \"We are applying for WSQ funding support for this new course Creating Stunning Digital Graphics with Canva according to Visual Communication DSN-COM-4005-1.1 under Design Framework.\"\n    New = \"We are applying for WSQ funding support for this new course AZ-900 Microsoft Azure Fundamentals Certification according to Cloud Computing ICT-DIT-4020-1.1 under Infocomm Technology Framework.\"\n  },\n  @{\n    Old = \"A significant challenge in the field is the inability to create coherent storyboards that clearly communicate visual intentions. Many designers lack the skills to evaluate and enhance their designs, leading to missed opportunities for improvement. This gap can hinder collaboration among team members and affect the overall quality of visual projects.\"\n    New = \"A prevalent issue in the industry is the inadequate integration of existing systems with cloud computing components. This often leads to compatibility problems and hinders the full utilization of cloud capabilities. Furthermore, organizations frequently lack the necessary skills to conduct impact analysis, which is crucial for understanding the implications of cloud adoption on their operations.\"\n  },\n  @{\n    Old = \"By emphasizing the creation of storyboards and task flows, this course addresses the need for clear visual communication. Participants will learn to evaluate their designs critically and suggest enhancements, fostering a collaborative environment that improves project outcomes.\"\n    New = \"The course equips participants with the skills to draft specifications that address performance and security requirements, ensuring that cloud solutions are tailored to organizational needs. Additionally, it emphasizes the importance of aligning existing systems with cloud components, which can significantly improve integration and overall system performance.\"\n  }\n)\n\nforeach ($para in $d.Paragraphs) {\n  $ptext = $para.Range.Text.TrimEnd([char]13, [char]7)\n  foreach ($rep in $replacements) {\n    if ($ptext -eq $rep.Old) {\n      $para.Range.Text = $rep.New\n      break\n    }\n  }\n}\n"}
